$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: swap the two list entries (logItem, description, man_hrs)
$ws.Range("B2").Value = "['HMV23/000059/0923/3', 'HMV23/000055/0923/17']"
$ws.Range("C2").Value = "['DURING ARRIVAL INSPECTION FOUND STATIC DISCHARGERS DAMAGED AT LOCATION : \n1) LH HORZ STAB TYPE A DISCHARGER QTY# 2 & TYPE B DISCHARGER QTY#3 \n2) RH HORZ STAB TYPE A DISCHARGER QTY#1 \n3)VERTICAL STABILISER TYPE A DISCHARGER QTY 01 & TYPE B DISCHARGER QTY 02.', 'DURING ARRIVAL INSPECTION FOUND STATIC DISCHARGERS DAMAGED AT LOCATION LH WING TIP TRAILING EDGE TYPE B DISCHARGER QTY#1, RH HORZ STAB TYPE A DISCHARGER QTY#1 TYPE B DISCHARGER QTY#1, LH HORZ STAB TYPE A DISCHARGER QTY#1 TYPE B DISCHARGER QTY#3, VERT STAB TYPE A DISCHARGER QTY#1, RUDDER TYPE B DISCHARGER QTY#1.']"
$ws.Range("D2").Value = "[5.0, 7.0]"

# Row 4: reorder 4 entries from [a,b,c,d] to [a,d,b,c]
$ws.Range("B4").Value = "['HMV23/000055/0923/2', 'HMV23/000043/0823/12', 'HMV23/000055/0923/6', 'HMV23/000059/0923/5']"
$ws.Range("C4").Value = "['DURING INSPECTION OBSERVED LH WING SLAT#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'WHILE ARRIVAL INSPECTION OBSERVED RH WING SLAT#5 OUTBOARD WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'DURING INSPECTION, OBSERVED RH SIDE WING SLAT NO#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.', 'DURING ARRIVAL INSPECTION , OBSERVED FOLLOWINS DICREPANCIES.\n\n1) RH SIDE WING NO#1 SLAT, NO.3 TRACK SQURE SEAL UPPER BULB SEAL DEGRADED.\n2) RH SIIDE WING SLAT NO#5 OUT BOARD WEATHER SEAL ERODED. SAME TO BE REPLACED.']"
$ws.Range("D4").Value = "[2.0, 2.0, 3.0, 0.0]"

# Row 6: swap the two list entries
$ws.Range("B6").Value = "['HMV23/000059/0923/7', 'HMV23/000055/0923/12']"
$ws.Range("C6").Value = "['DURING INSPECTION FOUND FOLLOWING AFT CARGO FLOOR PANEL DAMAGED , TO BE FABRICATED/REPLACED\n\n151LF P/N D5367521106800 \n\n151GF P/N D5367521100400 \n\n151DF P/N D5367521100200 \n\n151EF P/N D5367402100000', 'DURING INSPECTION FOUND FOLLOWING AFT CARGO FLOOR PANEL DAMAGED , TO BE FABRICATED/REPLACED\n\n151LF P/N D5367521106800\n\n151GF P/N D5367521100400\n\n151DF P/N D5367521100200\n\n151EF P/N D5367402100000\n\n151QF P/N D5367403500000']"

# Row 7: reorder 4 entries from [a,b,c,d] to [a,c,b,d] (logItem and man_hrs only; description unchanged)
$ws.Range("B7").Value = "['HMV23/000059/0923/8', 'HMV23/000055/0923/14', 'HMV23/000059/0923/9', 'HMV23/000055/0923/13']"
$ws.Range("D7").Value = "[5.0, 5.0, 7.0, 5.0]"

# Row 8: swap the two list entries
$ws.Range("B8").Value = "['HMV23/000055/0923/22', 'HMV23/000059/0923/10']"
$ws.Range("C8").Value = "['DURING OPP0RTUNITY INSPECTION FOUND BELLY PANEL HARDWARE DAMAGED/MISSING', 'DURING OPPORTUNITY INSPECTION FOUND BELLY PANEL HARDWARE DAMAGED/MISSING']"
$ws.Range("D8").Value = "[4.0, 2.0]"

# Row 11: swap the two list entries
$ws.Range("B11").Value = "['HMV23/000059/0923/15', 'HMV23/000055/0923/10']"
$ws.Range("C11").Value = "['DURING ARRIVAL CHECK FOUND EPSU BATTERIES QTY#2 FIN:10WL AND FIN: 11WL INOP.', 'DURING WEEKLY INSPECTION OBSERVED FIN: 11WL EPSU BATTERY INOP.']"
$ws.Range("D11").Value = "[6.0, 5.0]"

# Row 14: reorder 4 entries to match target order from diff
$ws.Range("B14").Value = "['HMV23/000055/0923/8', 'HMV23/000055/0923/1', 'HMV23/000055/0923/9', 'HMV23/000055/0923/7']"
$ws.Range("C14").Value = "['DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917', 'DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. ASSESSMENT TO BE CARRY OUT.\n\nSLAT NO #2 P/N- D5746092000500\n SA7205698', 'DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917', 'DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. DENT LOCATION, 3 INCH FROM T/E AND 67.5 INCH FROM OUTBD EDGE .ASSESSMENT TO BE CARRY OUT.\nSLAT NO #2 P/N- D5746092000500\n S/N- SA7205698']"
$ws.Range("D14").Value = "[3.0, 2.0, 7.0, 7.0]"

# Row 15: reorder 3 entries from [a,b,c] to [b,a,c]
$ws.Range("B15").Value = "['HMV23/000055/0923/3', 'HMV23/000055/0923/5', 'HMV23/000055/0923/4']"
$ws.Range("C15").Value = "['DURING ARRIVAL INSPECTION FOUND FOLLOWING WARNING MSG IN PFR:\n1. AUTO FLT AP OFF\nFURTHER RECTIFICATION TO BE CARRIED OUT.', 'DURING ARRIVAL INSPECTION FOUND FOLLOWING FAILURE MSG IN PFR:\n1. 38-31-41 TOILET ASSY LAV F\nFURTHER RECTIFICATION TO BE CARRIED OUT.', 'DURING ARRIVAL INSPECTION FOUND FOLLOWING FAILURE MSG IN PFR:\n1. 34-53-31 ADF1(2RP1)\nFURTHER RECTIFICATION TO BE CARRIED OUT.']"
$ws.Range("D15").Value = "[2.0, 4.0, 2.5]"

# Row 17: swap the two list entries
$ws.Range("B17").Value = "['HMV23/000043/0823/3', 'HMV23/000043/0823/5']"
$ws.Range("C17").Value = "['DURING ARRIVAL INSPECTION OBSERVED DELAMINATION ON FWD CARGO FLOOR PANELS 131AF,131EF,131FF,131JF.\n1)REMOVE FWD CARGO FLOOR PANEL 131AF,131EF,131FF,131JF FOR FABRICATION.\n2) INSTALL FLOOR PANEL 131AF,131EF,131FF,131JF POST FABRICATION.\n3) SEALS MISSING ON SOME FWD CARGO FLOOR PANELS.\n4) SEVERAL HARDWARE FOUND MISSING ON FWD CARGO SIDEWALL AND CEILING PANELS.', 'DURING ARRIVAL INSPECTION OBSERVED DELAMINATION ON AFT CARGO FLOOR PANELS 151SG,151DF,151GF, 51LF, 151PF.\n1)REMOVE AFT CARGO FLOOR PANELS 151SG,151DF,151GF, 51LF, 151PF FOR FABRICATION.\n2) INSTALL AFT CARGO FLOOR PANELS 151SG,151DF,151GF, 51LF, 151PF POST FABRICATION.\n3) SEVERAL HARDWARE FOUND MISSING ON AFT CARGO SIDEWALL AND CEILING PANELS.']"
